$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = -0.72
$wsSummary.Range("B6").Value = 25
$wsSummary.Range("B9").Value = 24

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 25
$wsStatus.Range("G4").Value = 24

# --- New trade row values (identical on both trade-log sheets) ---
$tradeNum      = 25
$tradeDate     = "2026-02-17"
$tradeTime     = "08:22:27"
$tradeStrategy = "MarketMaking"
$tradeSide     = "UP"
$entryPrice    = 0.68
$exitPrice     = 0.68
$tradeStatus   = "CLOSED"
$pnlPct        = 0
$pnlDollar     = 0
$capitalAfter  = 99.09999999999999
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = "early_exit"
$duration      = 0.13

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value  = $tradeNum
    # Copy the Date cell from the row above instead of assigning the string
    # directly - every prior row shares the same "2026-02-17" date text, and
    # a direct string assignment gets auto-parsed into a date serial number.
    $ws.Cells.Item($row - 1, 2).Copy($ws.Cells.Item($row, 2))
    $ws.Cells.Item($row, 3).Value  = $tradeTime
    $ws.Cells.Item($row, 4).Value  = $tradeStrategy
    $ws.Cells.Item($row, 5).Value  = $tradeSide
    $ws.Cells.Item($row, 6).Value  = $entryPrice
    $ws.Cells.Item($row, 7).Value  = $exitPrice
    $ws.Cells.Item($row, 8).Value  = $tradeStatus
    $ws.Cells.Item($row, 9).Value  = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades 26

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking 26
